$wb = $excel.ActiveWorkbook

# --- Sheet "Features" ---
$ws1 = $wb.Worksheets.Item("Features")

$ws1.Range("B2:E2").NumberFormat = "@"
$ws1.Range("B2").Value = "0,756"
$ws1.Range("C2").Value = "0,886"
$ws1.Range("D2").Value = "0,816"
$ws1.Range("E2").Value = "0,994"

$ws1.Range("B3:E3").NumberFormat = "@"
$ws1.Range("B3").Value = "0,800"
$ws1.Range("C3").Value = "0,914"
$ws1.Range("D3").Value = "0,853"
$ws1.Range("E3").Value = "0,889"

$ws1.Range("B4:E4").NumberFormat = "@"
$ws1.Range("B4").Value = "0,698"
$ws1.Range("C4").Value = "0,882"
$ws1.Range("D4").Value = "0,779"
$ws1.Range("E4").Value = "0,968"

$ws1.Range("B5:E5").NumberFormat = "@"
$ws1.Range("B5").Value = "0,750"
$ws1.Range("C5").Value = "0,882"
$ws1.Range("D5").Value = "0,811"
$ws1.Range("E5").Value = "0,873"

$ws1.Range("B6:E6").NumberFormat = "@"
$ws1.Range("B6").Value = "0,417"
$ws1.Range("C6").Value = "0,366"
$ws1.Range("D6").Value = "0,390"
$ws1.Range("E6").Value = "0,992"

$ws1.Range("B7:E7").NumberFormat = "@"
$ws1.Range("B7").Value = "0,556"
$ws1.Range("C7").Value = "0,513"
$ws1.Range("D7").Value = "0,533"
$ws1.Range("E7").Value = "1,000"

$ws1.Range("B8:E8").NumberFormat = "@"
$ws1.Range("B8").Value = "0,267"
$ws1.Range("C8").Value = "0,364"
$ws1.Range("D8").Value = "0,308"
$ws1.Range("E8").Value = "0,947"

$ws1.Range("B9:E9").NumberFormat = "@"
$ws1.Range("B9").Value = "0,600"
$ws1.Range("C9").Value = "0,714"
$ws1.Range("D9").Value = "0,652"
$ws1.Range("E9").Value = "1,000"

$ws1.Range("B10:E10").NumberFormat = "@"
$ws1.Range("B10").Value = "0,385"
$ws1.Range("C10").Value = "0,588"
$ws1.Range("D10").Value = "0,465"
$ws1.Range("E10").Value = "0,871"

$ws1.Range("B11:E11").NumberFormat = "@"
$ws1.Range("B11").Value = "0,556"
$ws1.Range("C11").Value = "0,294"
$ws1.Range("D11").Value = "0,385"
$ws1.Range("E11").Value = "0,989"

$ws1.Range("B12:E12").NumberFormat = "@"
$ws1.Range("B12").Value = "0,667"
$ws1.Range("C12").Value = "0,333"
$ws1.Range("D12").Value = "0,444"
$ws1.Range("E12").Value = "0,972"

$ws1.Range("B13:E13").NumberFormat = "@"
$ws1.Range("B13").Value = "0,500"
$ws1.Range("C13").Value = "0,294"
$ws1.Range("D13").Value = "0,370"
$ws1.Range("E13").Value = "0,995"

$ws1.Range("B14:E14").NumberFormat = "@"
$ws1.Range("B14").Value = "0,600"
$ws1.Range("C14").Value = "0,333"
$ws1.Range("D14").Value = "0,429"
$ws1.Range("E14").Value = "0,973"

$ws1.Range("B15:E15").NumberFormat = "@"
$ws1.Range("B15").Value = "0,579"
$ws1.Range("C15").Value = "0,524"
$ws1.Range("D15").Value = "0,550"
$ws1.Range("E15").Value = "0,773"

$ws1.Range("B16:E16").NumberFormat = "@"
$ws1.Range("B16").Value = "0,308"
$ws1.Range("C16").Value = "0,500"
$ws1.Range("D16").Value = "0,381"
$ws1.Range("E16").Value = "0,922"

$ws1.Range("B17:E17").NumberFormat = "@"
$ws1.Range("B17").Value = "0,722"
$ws1.Range("C17").Value = "0,812"
$ws1.Range("D17").Value = "0,765"
$ws1.Range("E17").Value = "0,973"

$ws1.Range("B18:E18").NumberFormat = "@"
$ws1.Range("B18").Value = "0,500"
$ws1.Range("C18").Value = "0,357"
$ws1.Range("D18").Value = "0,417"
$ws1.Range("E18").Value = "1,000"

$ws1.Range("B19:E19").NumberFormat = "@"
$ws1.Range("B19").Value = "0,478"
$ws1.Range("C19").Value = "0,733"
$ws1.Range("D19").Value = "0,579"
$ws1.Range("E19").Value = "0,843"

$ws1.Range("B20:E20").NumberFormat = "@"
$ws1.Range("B20").Value = "0,391"
$ws1.Range("C20").Value = "0,692"
$ws1.Range("D20").Value = "0,500"
$ws1.Range("E20").Value = "0,843"

$ws1.Range("B21:E21").NumberFormat = "@"
$ws1.Range("B21").Value = "0,750"
$ws1.Range("C21").Value = "0,429"
$ws1.Range("D21").Value = "0,545"
$ws1.Range("E21").Value = "0,848"

$ws1.Range("B22:E22").NumberFormat = "@"
$ws1.Range("B22").Value = "0,333"
$ws1.Range("C22").Value = "0,286"
$ws1.Range("D22").Value = "0,308"
$ws1.Range("E22").Value = "1,000"

$ws1.Range("B23:E23").NumberFormat = "@"
$ws1.Range("B23").Value = "0,500"
$ws1.Range("C23").Value = "0,400"
$ws1.Range("D23").Value = "0,444"
$ws1.Range("E23").Value = "0,685"

$ws1.Range("B24:E24").NumberFormat = "@"
$ws1.Range("B24").Value = "0,333"
$ws1.Range("C24").Value = "0,250"
$ws1.Range("D24").Value = "0,286"
$ws1.Range("E24").Value = "0,762"

$ws1.Range("B25:E25").NumberFormat = "@"
$ws1.Range("B25").Value = "0,500"
$ws1.Range("C25").Value = "0,400"
$ws1.Range("D25").Value = "0,444"
$ws1.Range("E25").Value = "0,685"

$ws1.Range("B26:E26").NumberFormat = "@"
$ws1.Range("B26").Value = "0,467"
$ws1.Range("C26").Value = "0,438"
$ws1.Range("D26").Value = "0,452"
$ws1.Range("E26").Value = "1,000"

$ws1.Range("B27:E27").NumberFormat = "@"
$ws1.Range("B27").Value = "0,400"
$ws1.Range("C27").Value = "0,222"
$ws1.Range("D27").Value = "0,286"
$ws1.Range("E27").Value = "1,000"

$ws1.Range("B28:E28").NumberFormat = "@"
$ws1.Range("B28").Value = "0,167"
$ws1.Range("C28").Value = "1,000"
$ws1.Range("D28").Value = "0,286"
$ws1.Range("E28").Value = "0,167"

$ws1.Range("B29:E29").NumberFormat = "@"
$ws1.Range("B29").Value = "0,083"
$ws1.Range("C29").Value = "0,125"
$ws1.Range("D29").Value = "0,100"
$ws1.Range("E29").Value = "1,000"

$ws1.Range("E30:E30").NumberFormat = "@"
$ws1.Range("E30").Value = "1,000"

$ws1.Range("E31:E31").NumberFormat = "@"
$ws1.Range("E31").Value = "0,884"

$ws1.Range("B32:E32").NumberFormat = "@"
$ws1.Range("B32").Value = "0,133"
$ws1.Range("C32").Value = "1,000"
$ws1.Range("D32").Value = "0,235"
$ws1.Range("E32").Value = "0,133"

$ws1.Range("B33:E33").NumberFormat = "@"
$ws1.Range("B33").Value = "0,500"
$ws1.Range("C33").Value = "1,000"
$ws1.Range("D33").Value = "0,667"
$ws1.Range("E33").Value = "0,500"

$ws1.Range("B34:E34").NumberFormat = "@"
$ws1.Range("B34").Value = "0,400"
$ws1.Range("C34").Value = "0,500"
$ws1.Range("D34").Value = "0,444"
$ws1.Range("E34").Value = "1,000"

$ws1.Range("B35:E35").NumberFormat = "@"
$ws1.Range("B35").Value = "0,200"
$ws1.Range("C35").Value = "0,200"
$ws1.Range("D35").Value = "0,200"
$ws1.Range("E35").Value = "1,000"

$ws1.Range("B36:E36").NumberFormat = "@"
$ws1.Range("B36").Value = "0,500"
$ws1.Range("C36").Value = "0,375"
$ws1.Range("D36").Value = "0,429"
$ws1.Range("E36").Value = "1,000"

$ws1.Range("B37:E37").NumberFormat = "@"
$ws1.Range("B37").Value = "0,500"
$ws1.Range("C37").Value = "0,143"
$ws1.Range("D37").Value = "0,222"
$ws1.Range("E37").Value = "1,000"

$ws1.Range("B38:E38").NumberFormat = "@"
$ws1.Range("B38").Value = "0,333"
$ws1.Range("C38").Value = "0,250"
$ws1.Range("D38").Value = "0,286"
$ws1.Range("E38").Value = "1,000"

$ws1.Range("B39:E39").NumberFormat = "@"
$ws1.Range("B39").Value = "0,500"
$ws1.Range("C39").Value = "0,667"
$ws1.Range("D39").Value = "0,571"
$ws1.Range("E39").Value = "1,000"

$ws1.Range("B40:E40").NumberFormat = "@"
$ws1.Range("B40").Value = "0,500"
$ws1.Range("C40").Value = "1,000"
$ws1.Range("D40").Value = "0,667"
$ws1.Range("E40").Value = "0,500"

$ws1.Range("B41:E41").NumberFormat = "@"
$ws1.Range("B41").Value = "0,750"
$ws1.Range("C41").Value = "1,000"
$ws1.Range("D41").Value = "0,857"
$ws1.Range("E41").Value = "0,750"

$ws1.Range("B42:E42").NumberFormat = "@"
$ws1.Range("B42").Value = "0,500"
$ws1.Range("C42").Value = "0,333"
$ws1.Range("D42").Value = "0,400"
$ws1.Range("E42").Value = "1,000"

$ws1.Range("B43:E43").NumberFormat = "@"
$ws1.Range("B43").Value = "0,667"
$ws1.Range("C43").Value = "1,000"
$ws1.Range("D43").Value = "0,800"
$ws1.Range("E43").Value = "0,667"

$ws1.Range("E45:E45").NumberFormat = "@"
$ws1.Range("E45").Value = "0,863"

$ws1.Range("B47:E47").NumberFormat = "@"
$ws1.Range("B47").Value = "0,500"
$ws1.Range("C47").Value = "0,250"
$ws1.Range("D47").Value = "0,333"
$ws1.Range("E47").Value = "1,000"

$ws1.Range("B48:E48").NumberFormat = "@"
$ws1.Range("B48").Value = "0,500"
$ws1.Range("C48").Value = "0,250"
$ws1.Range("D48").Value = "0,333"
$ws1.Range("E48").Value = "1,000"

$ws1.Range("E51:E51").NumberFormat = "@"
$ws1.Range("E51").Value = "0,131"

$ws1.Range("B55:E55").NumberFormat = "@"
$ws1.Range("B55").Value = "0,500"
$ws1.Range("C55").Value = "1,000"
$ws1.Range("D55").Value = "0,667"
$ws1.Range("E55").Value = "0,500"

$ws1.Range("B56:E56").NumberFormat = "@"
$ws1.Range("B56").Value = "0,500"
$ws1.Range("C56").Value = "0,500"
$ws1.Range("D56").Value = "0,500"
$ws1.Range("E56").Value = "0,452"

$ws1.Range("E68:E68").NumberFormat = "@"
$ws1.Range("E68").Value = "1,000"

# --- Sheet "Global Metrics" ---
$ws2 = $wb.Worksheets.Item("Global Metrics")

$ws2.Range("B2:E2").NumberFormat = "@"
$ws2.Range("B2").Value = "0,347"
$ws2.Range("C2").Value = "0,678"
$ws2.Range("D2").Value = "0,509"
$ws2.Range("E2").Value = "0,908"

